$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C3").Value = "Quản lý sản phẩm"
$ws.Range("C4").Select()
